$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = $null

$ws.Range("H70").Value = 3617.5293
$ws.Range("I70").Value = 2833.2222
$ws.Range("K70").Value = 8499.6666
$ws.Range("M70").Value = -8229.6666

$ws.Range("H73").Value = 3617.5293
$ws.Range("I73").Value = 2833.2222
$ws.Range("K73").Value = 8499.6666
$ws.Range("M73").Value = -7563.6666

$ws.Range("H118").Value = 329.92307
$ws.Range("I118").Value = 379.9091
$ws.Range("K118").Value = 1139.7273
$ws.Range("M118").Value = 517.2727

$ws.Range("H125").Value = 5421.7144
$ws.Range("I125").Value = 4325.5
$ws.Range("K125").Value = 38929.5
$ws.Range("M125").Value = -36469.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.6666
$ws.Range("I45").Value = 1999.6666
$ws.Range("K45").Value = 1999.6666
$ws.Range("M45").Value = -1622.6666

$ws.Range("H61").Value = 7011
$ws.Range("I61").Value = 7011
$ws.Range("K61").Value = 7011
$ws.Range("M61").Value = -6799

$ws.Range("H63").Value = 1802.5625
$ws.Range("I63").Value = 1413
$ws.Range("J63").Value = 2659.6
$ws.Range("K63").Value = 1413
$ws.Range("L63").Value = 2659.6
$ws.Range("M63").Value = -727
$ws.Range("N63").Value = -4031.6

$ws.Range("H66").Value = 1802.5625
$ws.Range("I66").Value = 1413
$ws.Range("J66").Value = 2659.6
$ws.Range("K66").Value = 7065
$ws.Range("L66").Value = 13298
$ws.Range("M66").Value = -3633
$ws.Range("N66").Value = -20162

$ws.Range("H74").Value = 5331.923
$ws.Range("J74").Value = 7500
$ws.Range("L74").Value = 7500
$ws.Range("N74").Value = -9248

$ws.Range("H77").Value = 5331.923
$ws.Range("J77").Value = 7500
$ws.Range("L77").Value = 37500
$ws.Range("N77").Value = -46236

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -9400

$ws.Range("H136").Value = 7011
$ws.Range("I136").Value = 7011
$ws.Range("K136").Value = 21033
$ws.Range("M136").Value = -18483

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5900
$ws.Range("I99").Value = 5900
$ws.Range("K99").Value = 5900
$ws.Range("M99").Value = -4402

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = $null

$ws.Range("H134").Value = 10158.375
$ws.Range("I134").Value = 8835.6
$ws.Range("K134").Value = 26506.8
$ws.Range("M134").Value = -23971.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4674.615
$ws.Range("I31").Value = 3882.8
$ws.Range("J31").Value = 5754.364
$ws.Range("K31").Value = 3882.8
$ws.Range("L31").Value = 5754.364
$ws.Range("M31").Value = -3587.8
$ws.Range("N31").Value = -6344.364

$ws.Range("H34").Value = 4674.615
$ws.Range("I34").Value = 3882.8
$ws.Range("J34").Value = 5754.364
$ws.Range("K34").Value = 3882.8
$ws.Range("L34").Value = 5754.364
$ws.Range("M34").Value = -3680.8
$ws.Range("N34").Value = -6158.364

$ws.Range("H58").Value = 1200
$ws.Range("I58").Value = 1200
$ws.Range("K58").Value = 1200
$ws.Range("M58").Value = -997

$ws.Range("H122").Value = 2275
$ws.Range("I122").Value = 2275
$ws.Range("K122").Value = 6825
$ws.Range("M122").Value = -4375

$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 1200
$ws.Range("K136").Value = 3600
$ws.Range("M136").Value = -1050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 97498.336
$ws.Range("J37").Value = 97498.336
$ws.Range("L37").Value = 292495.008
$ws.Range("N37").Value = -292719.008

$ws.Range("H50").Value = 616.3333
$ws.Range("I50").Value = 616.3333
$ws.Range("K50").Value = 1848.9999
$ws.Range("M50").Value = -1367.9999

$ws.Range("H53").Value = 616.3333
$ws.Range("I53").Value = 616.3333
$ws.Range("K53").Value = 1848.9999
$ws.Range("M53").Value = -1367.9999

$ws.Range("H116").Value = 1590.5555
$ws.Range("I116").Value = 720.5
$ws.Range("J116").Value = 3330.6667
$ws.Range("K116").Value = 2161.5
$ws.Range("L116").Value = 9992.000100000001
$ws.Range("M116").Value = 1280.5
$ws.Range("N116").Value = -16876.0001

$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3455.7827
$ws.Range("I46").Value = 3599.9
$ws.Range("J46").Value = 3344.923
$ws.Range("K46").Value = 3599.9
$ws.Range("L46").Value = 3344.923
$ws.Range("M46").Value = -3411.9
$ws.Range("N46").Value = -3720.923

$ws.Range("H122").Value = 2700.8
$ws.Range("I122").Value = 2700.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8102.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5652.400000000001
$ws.Range("N122").Value = $null

$ws.Range("H127").Value = 79497
$ws.Range("J127").Value = 79497
$ws.Range("L127").Value = 79497
$ws.Range("N127").Value = -89417

$ws.Range("H132").Value = 6636.385
$ws.Range("I132").Value = 4110.4287
$ws.Range("J132").Value = 9583.333000000001
$ws.Range("K132").Value = 12331.2861
$ws.Range("L132").Value = 28749.999
$ws.Range("M132").Value = -9801.286100000001
$ws.Range("N132").Value = -33809.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4289.8
$ws.Range("J62").Value = 3000
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 4289.8
$ws.Range("J65").Value = 3000
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null

$ws.Range("H101").Value = 59999
$ws.Range("J101").Value = 59999
$ws.Range("L101").Value = 59999
$ws.Range("N101").Value = -66489

$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50

$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1250
$ws.Range("K136").Value = 3750
$ws.Range("M136").Value = -1200
